$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value2 = 293.7577056666667
$ws.Range("H2").Value2 = 881.273117
$ws.Range("I2").Value2 = 0.9369756110667984
$ws.Range("J2").Value2 = 0.940660486426629
$ws.Range("M2").Value2 = 32.30926933333333
$ws.Range("N2").Value2 = 96.927808
$ws.Range("O2").Value2 = 0.4126751607889972
$ws.Range("P2").Value2 = 0.4140835826598944
$ws.Range("Q2").Value2 = 9491.096831126391
$ws.Range("R2").Value2 = 85419.87148013753
$ws.Range("S2").Value2 = 0.3866665609523599
$ws.Range("T2").Value2 = 0.3895120642861375

# Row 3
$ws.Range("G3").Value2 = 293.7577056666667
$ws.Range("H3").Value2 = 881.273117
$ws.Range("I3").Value2 = 0.9369756110667984
$ws.Range("J3").Value2 = 0.940660486426629
$ws.Range("O3").Value2 = 0.03337290046497914
$ws.Range("P3").Value2 = 0.03348679906459509
$ws.Range("Q3").Value2 = 767.541785754858
$ws.Range("R3").Value2 = 6907.876071793722
$ws.Range("S3").Value2 = 0.03126959380624527
$ws.Range("T3").Value2 = 0.03149970869697281

# Row 4
$ws.Range("G4").Value2 = 293.7577056666667
$ws.Range("H4").Value2 = 881.273117
$ws.Range("I4").Value2 = 0.9369756110667984
$ws.Range("J4").Value2 = 0.940660486426629
$ws.Range("M4").Value2 = 12.89411533333333
$ws.Range("N4").Value2 = 38.682346
$ws.Range("O4").Value2 = 0.1646920907903501
$ws.Range("P4").Value2 = 0.1652541695502867
$ws.Range("Q4").Value2 = 3787.745736921386
$ws.Range("R4").Value2 = 34089.71163229248
$ws.Range("S4").Value2 = 0.1543124724061569
$ws.Range("T4").Value2 = 0.1554480675132013

# Row 5
$ws.Range("G5").Value2 = 293.7577056666667
$ws.Range("H5").Value2 = 881.273117
$ws.Range("I5").Value2 = 0.9369756110667984
$ws.Range("J5").Value2 = 0.940660486426629
$ws.Range("M5").Value2 = 0.798886
$ws.Range("N5").Value2 = 1.597772
$ws.Range("O5").Value2 = 0.01020389551681842
$ws.Range("P5").Value2 = 0.006825813640948785
$ws.Range("Q5").Value2 = 234.6789184492206
$ws.Range("R5").Value2 = 1408.073510695324
$ws.Range("S5").Value2 = 0.009560801237132703
$ws.Range("T5").Value2 = 0.006420773179752403

# Row 6
$ws.Range("G6").Value2 = 293.7577056666667
$ws.Range("H6").Value2 = 881.273117
$ws.Range("I6").Value2 = 0.9369756110667984
$ws.Range("J6").Value2 = 0.940660486426629
$ws.Range("M6").Value2 = 29.67714566666667
$ws.Range("N6").Value2 = 89.03143700000001
$ws.Range("O6").Value2 = 0.3790559524388551
$ws.Range("P6").Value2 = 0.3803496350842752
$ws.Range("Q6").Value2 = 8717.89022177546
$ws.Range("R6").Value2 = 78461.01199597913
$ws.Range("S6").Value2 = 0.3551661826649035
$ws.Range("T6").Value2 = 0.3577798727505651

# Row 7
$ws.Range("I7").Value2 = 0.01525159481997056
$ws.Range("J7").Value2 = 0.01531157527761154
$ws.Range("M7").Value2 = 32.30926933333333
$ws.Range("N7").Value2 = 96.927808
$ws.Range("O7").Value2 = 0.4126751607889972
$ws.Range("P7").Value2 = 0.4140835826598944
$ws.Range("Q7").Value2 = 154.4910684501547
$ws.Range("R7").Value2 = 1390.419616051392
$ws.Range("S7").Value2 = 0.006293954344619986
$ws.Range("T7").Value2 = 0.006340271947120055

# Row 8
$ws.Range("I8").Value2 = 0.01525159481997056
$ws.Range("J8").Value2 = 0.01531157527761154
$ws.Range("O8").Value2 = 0.03337290046497914
$ws.Range("P8").Value2 = 0.03348679906459509
$ws.Range("S8").Value2 = 0.0005089899558590689
$ws.Range("T8").Value2 = 0.0005127356446837996

# Row 9
$ws.Range("I9").Value2 = 0.01525159481997056
$ws.Range("J9").Value2 = 0.01531157527761154
$ws.Range("M9").Value2 = 12.89411533333333
$ws.Range("N9").Value2 = 38.682346
$ws.Range("O9").Value2 = 0.1646920907903501
$ws.Range("P9").Value2 = 0.1652541695502867
$ws.Range("Q9").Value2 = 61.65492738367266
$ws.Range("R9").Value2 = 554.8943464530539
$ws.Range("S9").Value2 = 0.002511817038788224
$ws.Range("T9").Value2 = 0.002530301657008396

# Row 10
$ws.Range("I10").Value2 = 0.01525159481997056
$ws.Range("J10").Value2 = 0.01531157527761154
$ws.Range("M10").Value2 = 0.798886
$ws.Range("N10").Value2 = 1.597772
$ws.Range("O10").Value2 = 0.01020389551681842
$ws.Range("P10").Value2 = 0.006825813640948785
$ws.Range("Q10").Value2 = 3.819979660838
$ws.Range("R10").Value2 = 22.919877965028
$ws.Range("S10").Value2 = 0.0001556256800078286
$ws.Range("T10").Value2 = 0.000104513959394335

# Row 11
$ws.Range("I11").Value2 = 0.01525159481997056
$ws.Range("J11").Value2 = 0.01531157527761154
$ws.Range("M11").Value2 = 29.67714566666667
$ws.Range("N11").Value2 = 89.03143700000001
$ws.Range("O11").Value2 = 0.3790559524388551
$ws.Range("P11").Value2 = 0.3803496350842752
$ws.Range("Q11").Value2 = 141.9052190655404
$ws.Range("R11").Value2 = 1277.146971589863
$ws.Range("S11").Value2 = 0.005781207800695448
$ws.Range("T11").Value2 = 0.00582375206940496

# Row 12
$ws.Range("G12").Value2 = 7.906212666666666
$ws.Range("H12").Value2 = 23.718638
$ws.Range("I12").Value2 = 0.02521781829607561
$ws.Range("J12").Value2 = 0.02531699325449539
$ws.Range("M12").Value2 = 32.30926933333333
$ws.Range("N12").Value2 = 96.927808
$ws.Range("O12").Value2 = 0.4126751607889972
$ws.Range("P12").Value2 = 0.4140835826598944
$ws.Range("Q12").Value2 = 255.4439544539449
$ws.Range("R12").Value2 = 2298.995590085504
$ws.Range("S12").Value2 = 0.01040676722008072
$ws.Range("T12").Value2 = 0.01048335126899783

# Row 13
$ws.Range("G13").Value2 = 7.906212666666666
$ws.Range("H13").Value2 = 23.718638
$ws.Range("I13").Value2 = 0.02521781829607561
$ws.Range("J13").Value2 = 0.02531699325449539
$ws.Range("O13").Value2 = 0.03337290046497914
$ws.Range("P13").Value2 = 0.03348679906459509
$ws.Range("Q13").Value2 = 20.65766606856911
$ws.Range("R13").Value2 = 185.918994617122
$ws.Range("S13").Value2 = 0.0008415917399388613
$ws.Range("T13").Value2 = 0.0008477850660329964

# Row 14
$ws.Range("G14").Value2 = 7.906212666666666
$ws.Range("H14").Value2 = 23.718638
$ws.Range("I14").Value2 = 0.02521781829607561
$ws.Range("J14").Value2 = 0.02531699325449539
$ws.Range("M14").Value2 = 12.89411533333333
$ws.Range("N14").Value2 = 38.682346
$ws.Range("O14").Value2 = 0.1646920907903501
$ws.Range("P14").Value2 = 0.1652541695502867
$ws.Range("Q14").Value2 = 101.9436179738609
$ws.Range("R14").Value2 = 917.4925617647478
$ws.Range("S14").Value2 = 0.004153175220351836
$ws.Range("T14").Value2 = 0.004183738695781846

# Row 15
$ws.Range("G15").Value2 = 7.906212666666666
$ws.Range("H15").Value2 = 23.718638
$ws.Range("I15").Value2 = 0.02521781829607561
$ws.Range("J15").Value2 = 0.02531699325449539
$ws.Range("M15").Value2 = 0.798886
$ws.Range("N15").Value2 = 1.597772
$ws.Range("O15").Value2 = 0.01020389551681842
$ws.Range("P15").Value2 = 0.006825813640948785
$ws.Range("Q15").Value2 = 6.316162612422666
$ws.Range("R15").Value2 = 37.89697567453599
$ws.Range("S15").Value2 = 0.0002573199830552674
$ws.Range("T15").Value2 = 0.000172809077904343

# Row 16
$ws.Range("G16").Value2 = 7.906212666666666
$ws.Range("H16").Value2 = 23.718638
$ws.Range("I16").Value2 = 0.02521781829607561
$ws.Range("J16").Value2 = 0.02531699325449539
$ws.Range("M16").Value2 = 29.67714566666667
$ws.Range("N16").Value2 = 89.03143700000001
$ws.Range("O16").Value2 = 0.3790559524388551
$ws.Range("P16").Value2 = 0.3803496350842752
$ws.Range("Q16").Value2 = 234.6338249803118
$ws.Range("R16").Value2 = 2111.704424822806
$ws.Range("S16").Value2 = 0.009558964132648926
$ws.Range("T16").Value2 = 0.009629309145778378

# Row 17
$ws.Range("G17").Value2 = 3.6844455
$ws.Range("H17").Value2 = 7.368891000000001
$ws.Range("I17").Value2 = 0.01175198303639443
$ws.Range("J17").Value2 = 0.007865466969060864
$ws.Range("M17").Value2 = 32.30926933333333
$ws.Range("N17").Value2 = 96.927808
$ws.Range("O17").Value2 = 0.4126751607889972
$ws.Range("P17").Value2 = 0.4140835826598944
$ws.Range("Q17").Value2 = 119.041742003488
$ws.Range("R17").Value2 = 714.250452020928
$ws.Range("S17").Value2 = 0.004849751489133641
$ws.Range("T17").Value2 = 0.003256960741841783

# Row 18
$ws.Range("G18").Value2 = 3.6844455
$ws.Range("H18").Value2 = 7.368891000000001
$ws.Range("I18").Value2 = 0.01175198303639443
$ws.Range("J18").Value2 = 0.007865466969060864
$ws.Range("O18").Value2 = 0.03337290046497914
$ws.Range("P18").Value2 = 0.03348679906459509
$ws.Range("Q18").Value2 = 9.6268653520715
$ws.Range("R18").Value2 = 57.761192112429
$ws.Range("S18").Value2 = 0.0003921977601397148
$ws.Range("T18").Value2 = 0.0002633893119421509

# Row 19
$ws.Range("G19").Value2 = 3.6844455
$ws.Range("H19").Value2 = 7.368891000000001
$ws.Range("I19").Value2 = 0.01175198303639443
$ws.Range("J19").Value2 = 0.007865466969060864
$ws.Range("M19").Value2 = 12.89411533333333
$ws.Range("N19").Value2 = 38.682346
$ws.Range("O19").Value2 = 0.1646920907903501
$ws.Range("P19").Value2 = 0.1652541695502867
$ws.Range("Q19").Value2 = 47.507665216381
$ws.Range("R19").Value2 = 285.045991298286
$ws.Range("S19").Value2 = 0.001935458657196526
$ws.Range("T19").Value2 = 0.001299801212097363

# Row 20
$ws.Range("G20").Value2 = 3.6844455
$ws.Range("H20").Value2 = 7.368891000000001
$ws.Range("I20").Value2 = 0.01175198303639443
$ws.Range("J20").Value2 = 0.007865466969060864
$ws.Range("M20").Value2 = 0.798886
$ws.Range("N20").Value2 = 1.597772
$ws.Range("O20").Value2 = 0.01020389551681842
$ws.Range("P20").Value2 = 0.006825813640948785
$ws.Range("Q20").Value2 = 2.943451927713
$ws.Range("R20").Value2 = 11.773807710852
$ws.Range("S20").Value2 = 0.0001199160070187913
$ws.Range("T20").Value2 = 0.00005368821172984774

# Row 21
$ws.Range("G21").Value2 = 3.6844455
$ws.Range("H21").Value2 = 7.368891000000001
$ws.Range("I21").Value2 = 0.01175198303639443
$ws.Range("J21").Value2 = 0.007865466969060864
$ws.Range("M21").Value2 = 29.67714566666667
$ws.Range("N21").Value2 = 89.03143700000001
$ws.Range("O21").Value2 = 0.3790559524388551
$ws.Range("P21").Value2 = 0.3803496350842752
$ws.Range("Q21").Value2 = 109.3438258043945
$ws.Range("R21").Value2 = 656.0629548263671
$ws.Range("S21").Value2 = 0.004454659122905761
$ws.Range("T21").Value2 = 0.00299162749144972

# Row 22
$ws.Range("G22").Value2 = 3.386921
$ws.Range("H22").Value2 = 10.160763
$ws.Range("I22").Value2 = 0.01080299278076119
$ws.Range("J22").Value2 = 0.01084547807220323
$ws.Range("M22").Value2 = 32.30926933333333
$ws.Range("N22").Value2 = 96.927808
$ws.Range("O22").Value2 = 0.4126751607889972
$ws.Range("P22").Value2 = 0.4140835826598944
$ws.Range("Q22").Value2 = 109.4289427997226
$ws.Range("R22").Value2 = 984.860485197504
$ws.Range("S22").Value2 = 0.004458126782803002
$ws.Range("T22").Value2 = 0.00449093441579724

# Row 23
$ws.Range("G23").Value2 = 3.386921
$ws.Range("H23").Value2 = 10.160763
$ws.Range("I23").Value2 = 0.01080299278076119
$ws.Range("J23").Value2 = 0.01084547807220323
$ws.Range("O23").Value2 = 0.03337290046497914
$ws.Range("P23").Value2 = 0.03348679906459509
$ws.Range("Q23").Value2 = 8.849481536666332
$ws.Range("R23").Value2 = 79.64533382999699
$ws.Range("S23").Value2 = 0.0003605272027962316
$ws.Range("T23").Value2 = 0.0003631803449633418

# Row 24
$ws.Range("G24").Value2 = 3.386921
$ws.Range("H24").Value2 = 10.160763
$ws.Range("I24").Value2 = 0.01080299278076119
$ws.Range("J24").Value2 = 0.01084547807220323
$ws.Range("M24").Value2 = 12.89411533333333
$ws.Range("N24").Value2 = 38.682346
$ws.Range("O24").Value2 = 0.1646920907903501
$ws.Range("P24").Value2 = 0.1652541695502867
$ws.Range("Q24").Value2 = 43.67134999888866
$ws.Range("R24").Value2 = 393.0421499899979
$ws.Range("S24").Value2 = 0.001779167467856619
$ws.Range("T24").Value2 = 0.001792260472197789

# Row 25
$ws.Range("G25").Value2 = 3.386921
$ws.Range("H25").Value2 = 10.160763
$ws.Range("I25").Value2 = 0.01080299278076119
$ws.Range("J25").Value2 = 0.01084547807220323
$ws.Range("M25").Value2 = 0.798886
$ws.Range("N25").Value2 = 1.597772
$ws.Range("O25").Value2 = 0.01020389551681842
$ws.Range("P25").Value2 = 0.006825813640948785
$ws.Range("Q25").Value2 = 2.705763770006
$ws.Range("R25").Value2 = 16.234582620036
$ws.Range("S25").Value2 = 0.0001102326096038309
$ws.Range("T25").Value2 = 0.00007402921216785576

# Row 26
$ws.Range("G26").Value2 = 3.386921
$ws.Range("H26").Value2 = 10.160763
$ws.Range("I26").Value2 = 0.01080299278076119
$ws.Range("J26").Value2 = 0.01084547807220323
$ws.Range("M26").Value2 = 29.67714566666667
$ws.Range("N26").Value2 = 89.03143700000001
$ws.Range("O26").Value2 = 0.3790559524388551
$ws.Range("P26").Value2 = 0.3803496350842752
$ws.Range("Q26").Value2 = 100.5141478784923
$ws.Range("R26").Value2 = 904.6273309064311
$ws.Range("S26").Value2 = 0.00409493871770151
$ws.Range("T26").Value2 = 0.003577972480727868
